$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 17-19 (old Resolving-Mac block removed)
$ws.Rows.Item(17).Resize(3).Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il34"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"2"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"4.311203"
$ws.Range("H2").Value = [double]"8.622406"
$ws.Range("I2").Value = [double]"0.168599917917489"
$ws.Range("J2").Value = [double]"0.162674471805938"
$ws.Range("K2").Value = [double]"1"
$ws.Range("L2").Value = [double]"0.5"
$ws.Range("M2").Value = [double]"0.07111099999999999"
$ws.Range("N2").Value = [double]"0.142222"
$ws.Range("O2").Value = [double]"0.02711460746047303"
$ws.Range("P2").Value = [double]"0.02678527708115022"
$ws.Range("Q2").Value = [double]"0.306573956533"
$ws.Range("R2").Value = [double]"1.226295826132"
$ws.Range("S2").Value = [double]"0.004571520592200689"
$ws.Range("T2").Value = [double]"0.004357280801351809"

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il34"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"2"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"4.311203"
$ws.Range("H3").Value = [double]"8.622406"
$ws.Range("I3").Value = [double]"0.168599917917489"
$ws.Range("J3").Value = [double]"0.162674471805938"
$ws.Range("K3").Value = [double]"2"
$ws.Range("L3").Value = [double]"0.6666666666666666"
$ws.Range("M3").Value = [double]"0.06449100000000001"
$ws.Range("N3").Value = [double]"0.193473"
$ws.Range("O3").Value = [double]"0.0245904030281302"
$ws.Range("P3").Value = [double]"0.03643759694506741"
$ws.Range("Q3").Value = [double]"0.278033792673"
$ws.Range("R3").Value = [double]"1.668202756038"
$ws.Range("S3").Value = [double]"0.004145939932100725"
$ws.Range("T3").Value = [double]"0.005927466836916501"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il34"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = [double]"2"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"4.311203"
$ws.Range("H4").Value = [double]"8.622406"
$ws.Range("I4").Value = [double]"0.168599917917489"
$ws.Range("J4").Value = [double]"0.162674471805938"
$ws.Range("K4").Value = [double]"2"
$ws.Range("L4").Value = [double]"1"
$ws.Range("M4").Value = [double]"2.4870065"
$ws.Range("N4").Value = [double]"4.974013"
$ws.Range("O4").Value = [double]"0.9482949895113968"
$ws.Range("P4").Value = [double]"0.9367771259737823"
$ws.Range("Q4").Value = [double]"10.7219898838195"
$ws.Range("R4").Value = [double]"42.887959535278"
$ws.Range("S4").Value = [double]"0.1598824573931876"
$ws.Range("T4").Value = [double]"0.1523897241676697"

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il34"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"1.560702"
$ws.Range("H5").Value = [double]"4.682106"
$ws.Range("I5").Value = [double]"0.0610349893274942"
$ws.Range("J5").Value = [double]"0.08833487085732371"
$ws.Range("K5").Value = [double]"1"
$ws.Range("L5").Value = [double]"0.5"
$ws.Range("M5").Value = [double]"0.07111099999999999"
$ws.Range("N5").Value = [double]"0.142222"
$ws.Range("O5").Value = [double]"0.02711460746047303"
$ws.Range("P5").Value = [double]"0.02678527708115022"
$ws.Range("Q5").Value = [double]"0.110983079922"
$ws.Range("R5").Value = [double]"0.665898479532"
$ws.Range("S5").Value = [double]"0.001654939776969166"
$ws.Range("T5").Value = [double]"0.002366073991841037"

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il34"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"1.560702"
$ws.Range("H6").Value = [double]"4.682106"
$ws.Range("I6").Value = [double]"0.0610349893274942"
$ws.Range("J6").Value = [double]"0.08833487085732371"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.06449100000000001"
$ws.Range("N6").Value = [double]"0.193473"
$ws.Range("O6").Value = [double]"0.0245904030281302"
$ws.Range("P6").Value = [double]"0.03643759694506741"
$ws.Range("Q6").Value = [double]"0.100651232682"
$ws.Range("R6").Value = [double]"0.9058610941380001"
$ws.Range("S6").Value = [double]"0.001500874986380708"
$ws.Range("T6").Value = [double]"0.003218710420493742"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il34"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = [double]"3"
$ws.Range("F7").Value = [double]"1"
$ws.Range("G7").Value = [double]"1.560702"
$ws.Range("H7").Value = [double]"4.682106"
$ws.Range("I7").Value = [double]"0.0610349893274942"
$ws.Range("J7").Value = [double]"0.08833487085732371"
$ws.Range("K7").Value = [double]"2"
$ws.Range("L7").Value = [double]"1"
$ws.Range("M7").Value = [double]"2.4870065"
$ws.Range("N7").Value = [double]"4.974013"
$ws.Range("O7").Value = [double]"0.9482949895113968"
$ws.Range("P7").Value = [double]"0.9367771259737823"
$ws.Range("Q7").Value = [double]"3.881476018563"
$ws.Range("R7").Value = [double]"23.288856111378"
$ws.Range("S7").Value = [double]"0.05787917456414433"
$ws.Range("T7").Value = [double]"0.08275008644498892"

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Il34"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = [double]"2"
$ws.Range("F8").Value = [double]"1"
$ws.Range("G8").Value = [double]"19.3965805"
$ws.Range("H8").Value = [double]"38.793161"
$ws.Range("I8").Value = [double]"0.7585497319843134"
$ws.Range("J8").Value = [double]"0.7318904926719657"
$ws.Range("K8").Value = [double]"1"
$ws.Range("L8").Value = [double]"0.5"
$ws.Range("M8").Value = [double]"0.07111099999999999"
$ws.Range("N8").Value = [double]"0.142222"
$ws.Range("O8").Value = [double]"0.02711460746047303"
$ws.Range("P8").Value = [double]"0.02678527708115022"
$ws.Range("Q8").Value = [double]"1.3793102359355"
$ws.Range("R8").Value = [double]"5.517240943741999"
$ws.Range("S8").Value = [double]"0.02056777822200169"
$ws.Range("T8").Value = [double]"0.01960388963927815"

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Il34"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = [double]"2"
$ws.Range("F9").Value = [double]"1"
$ws.Range("G9").Value = [double]"19.3965805"
$ws.Range("H9").Value = [double]"38.793161"
$ws.Range("I9").Value = [double]"0.7585497319843134"
$ws.Range("J9").Value = [double]"0.7318904926719657"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"0.06449100000000001"
$ws.Range("N9").Value = [double]"0.193473"
$ws.Range("O9").Value = [double]"0.0245904030281302"
$ws.Range("P9").Value = [double]"0.03643759694506741"
$ws.Range("Q9").Value = [double]"1.2509048730255"
$ws.Range("R9").Value = [double]"7.505429238153"
$ws.Range("S9").Value = [double]"0.01865304362637441"
$ws.Range("T9").Value = [double]"0.0266683307799079"

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Il34"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = [double]"2"
$ws.Range("F10").Value = [double]"1"
$ws.Range("G10").Value = [double]"19.3965805"
$ws.Range("H10").Value = [double]"38.793161"
$ws.Range("I10").Value = [double]"0.7585497319843134"
$ws.Range("J10").Value = [double]"0.7318904926719657"
$ws.Range("K10").Value = [double]"2"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"2.4870065"
$ws.Range("N10").Value = [double]"4.974013"
$ws.Range("O10").Value = [double]"0.9482949895113968"
$ws.Range("P10").Value = [double]"0.9367771259737823"
$ws.Range("Q10").Value = [double]"48.23942178127325"
$ws.Range("R10").Value = [double]"192.957687125093"
$ws.Range("S10").Value = [double]"0.7193289101359374"
$ws.Range("T10").Value = [double]"0.6856182722527796"

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "Il34"
$ws.Range("C11").Value = "Ptprz1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = [double]"3"
$ws.Range("F11").Value = [double]"1"
$ws.Range("G11").Value = [double]"0.220446"
$ws.Range("H11").Value = [double]"0.661338"
$ws.Range("I11").Value = [double]"0.008621068760909376"
$ws.Range("J11").Value = [double]"0.0124771217958416"
$ws.Range("K11").Value = [double]"1"
$ws.Range("L11").Value = [double]"0.5"
$ws.Range("M11").Value = [double]"0.07111099999999999"
$ws.Range("N11").Value = [double]"0.142222"
$ws.Range("O11").Value = [double]"0.02711460746047303"
$ws.Range("P11").Value = [double]"0.02678527708115022"
$ws.Range("Q11").Value = [double]"0.015676135506"
$ws.Range("R11").Value = [double]"0.09405681303599998"
$ws.Range("S11").Value = [double]"0.0002337568953418044"
$ws.Range("T11").Value = [double]"0.000334203164476876"

# Row 12
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "Il34"
$ws.Range("C12").Value = "Ptprz1"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = [double]"3"
$ws.Range("F12").Value = [double]"1"
$ws.Range("G12").Value = [double]"0.220446"
$ws.Range("H12").Value = [double]"0.661338"
$ws.Range("I12").Value = [double]"0.008621068760909376"
$ws.Range("J12").Value = [double]"0.0124771217958416"
$ws.Range("K12").Value = [double]"2"
$ws.Range("L12").Value = [double]"0.6666666666666666"
$ws.Range("M12").Value = [double]"0.06449100000000001"
$ws.Range("N12").Value = [double]"0.193473"
$ws.Range("O12").Value = [double]"0.0245904030281302"
$ws.Range("P12").Value = [double]"0.03643759694506741"
$ws.Range("Q12").Value = [double]"0.014216782986"
$ws.Range("R12").Value = [double]"0.127951046874"
$ws.Range("S12").Value = [double]"0.0002119955553639846"
$ws.Range("T12").Value = [double]"0.000454636335031392"

# Row 13
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "Il34"
$ws.Range("C13").Value = "Ptprz1"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = [double]"3"
$ws.Range("F13").Value = [double]"1"
$ws.Range("G13").Value = [double]"0.220446"
$ws.Range("H13").Value = [double]"0.661338"
$ws.Range("I13").Value = [double]"0.008621068760909376"
$ws.Range("J13").Value = [double]"0.0124771217958416"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"2.4870065"
$ws.Range("N13").Value = [double]"4.974013"
$ws.Range("O13").Value = [double]"0.9482949895113968"
$ws.Range("P13").Value = [double]"0.9367771259737823"
$ws.Range("Q13").Value = [double]"0.5482506348990001"
$ws.Range("R13").Value = [double]"3.289503809394"
$ws.Range("S13").Value = [double]"0.008175316310203587"
$ws.Range("T13").Value = [double]"0.01168828229633333"

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Il34"
$ws.Range("C14").Value = "Ptprz1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.08168"
$ws.Range("H14").Value = [double]"0.24504"
$ws.Range("I14").Value = [double]"0.003194292009794135"
$ws.Range("J14").Value = [double]"0.004623042868930905"
$ws.Range("K14").Value = [double]"1"
$ws.Range("L14").Value = [double]"0.5"
$ws.Range("M14").Value = [double]"0.07111099999999999"
$ws.Range("N14").Value = [double]"0.142222"
$ws.Range("O14").Value = [double]"0.02711460746047303"
$ws.Range("P14").Value = [double]"0.02678527708115022"
$ws.Range("Q14").Value = [double]"0.00580834648"
$ws.Range("R14").Value = [double]"0.03485007888"
$ws.Range("S14").Value = [double]"8.661197395969344E-05"
$ws.Range("T14").Value = [double]"0.0001238294842023499"

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Il34"
$ws.Range("C15").Value = "Ptprz1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.08168"
$ws.Range("H15").Value = [double]"0.24504"
$ws.Range("I15").Value = [double]"0.003194292009794135"
$ws.Range("J15").Value = [double]"0.004623042868930905"
$ws.Range("K15").Value = [double]"2"
$ws.Range("L15").Value = [double]"0.6666666666666666"
$ws.Range("M15").Value = [double]"0.06449100000000001"
$ws.Range("N15").Value = [double]"0.193473"
$ws.Range("O15").Value = [double]"0.0245904030281302"
$ws.Range("P15").Value = [double]"0.03643759694506741"
$ws.Range("Q15").Value = [double]"0.005267624880000001"
$ws.Range("R15").Value = [double]"0.04740862392"
$ws.Range("S15").Value = [double]"7.854892791037379E-05"
$ws.Range("T15").Value = [double]"0.0001684525727178724"

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Il34"
$ws.Range("C16").Value = "Ptprz1"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.08168"
$ws.Range("H16").Value = [double]"0.24504"
$ws.Range("I16").Value = [double]"0.003194292009794135"
$ws.Range("J16").Value = [double]"0.004623042868930905"
$ws.Range("K16").Value = [double]"2"
$ws.Range("L16").Value = [double]"1"
$ws.Range("M16").Value = [double]"2.4870065"
$ws.Range("N16").Value = [double]"4.974013"
$ws.Range("O16").Value = [double]"0.9482949895113968"
$ws.Range("P16").Value = [double]"0.9367771259737823"
$ws.Range("Q16").Value = [double]"0.20313869092"
$ws.Range("R16").Value = [double]"1.21883214552"
$ws.Range("S16").Value = [double]"0.003029131107924068"
$ws.Range("T16").Value = [double]"0.004330760812010682"
